# Vdd readout routine added
#
# 1. Typography sheet: the "Default" typography (row 4) now gets the same
#    Wildcard Characters / Wildcard Ranges as "Large"/"Huge" so numeric
#    Vdd readouts render correctly.
# 2. Translation sheet: drop the stale "SingleUseId3"/"Red" row (row 4),
#    which shifts every following row up by one, and append the two new
#    text entries needed for the Vdd readout: a "Vdd" label and a
#    "<value>V" value+unit text.

$wb = $excel.ActiveWorkbook

$typography = $wb.Worksheets.Item("Typography")
$typography.Range("G4").Value = "."
$typography.Range("I4").Value = "0-9"

$translation = $wb.Worksheets.Item("Translation")
$translation.Rows.Item(4).Delete()

$translation.Range("B21").Value = "SingleUseId25"
$translation.Range("C21").Value = "Default"
$translation.Range("D21").Value = "Right"
$translation.Range("E21").Value = "LTR"
$translation.Range("F21").Value = "Vdd"

$translation.Range("B22").Value = "SingleUseId26"
$translation.Range("C22").Value = "Default"
$translation.Range("D22").Value = "Right"
$translation.Range("E22").Value = "LTR"
$translation.Range("F22").Value = "<value>V"
